# Insert a new daily price record at row 420 (the sheet's data is stored
# most-recent-first, so inserting at the top of the data block represents a
# new day's observation). All subsequent rows shift down by one, which also
# grows the sheet's used range from A1:T522 to A1:T523.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A420:T420").EntireRow.Insert()

$ws.Cells.Item(420, 1).Value  = 9
$ws.Cells.Item(420, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(420, 3).Value  = "Metropolitana"
$ws.Cells.Item(420, 4).Value  = 44855
$ws.Cells.Item(420, 5).Value  = 13
$ws.Cells.Item(420, 6).Value  = "Fruta"
$ws.Cells.Item(420, 7).Value  = 100108
$ws.Cells.Item(420, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(420, 9).Value  = 100108002
$ws.Cells.Item(420, 10).Value = "Mango"
$ws.Cells.Item(420, 11).Value = "Sin especificar"
$ws.Cells.Item(420, 12).Value = "Primera"
$ws.Cells.Item(420, 13).Value = 450
$ws.Cells.Item(420, 14).Value = 6500
$ws.Cells.Item(420, 15).Value = 7000
$ws.Cells.Item(420, 16).Value = 6722
$ws.Cells.Item(420, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(420, 18).Value = "Brasil"
$ws.Cells.Item(420, 19).Value = 1680
$ws.Cells.Item(420, 20).Value = 4
